$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Pretax Income" label to "Pretax Income (Loss)" before deleting columns
$ws.Range("P1").Value = "Pretax Income (Loss)"

# Delete columns: Selling General & Admin (J), Other Operating Expenses (K),
# Other Expense / Income (O), Income Tax (Q)
$ws.Range("Q1:Q13").EntireColumn.Delete()
$ws.Range("O1:O13").EntireColumn.Delete()
$ws.Range("J1:K13").EntireColumn.Delete()

# Adjust view: scroll so column B is the top-left visible column, select M1
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("M1").Select()
